# Fruta / hortaliza, semanal
# Insert 3 new weekly report rows at the top of the data block (row 274),
# pushing the existing rows (274:301) down to (277:304).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 274.
$ws.Rows("274:276").Insert()

# Constant values shared by every data row in this sheet.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112021
$categoria = "Ají"
$unidad    = "`$/caja 25 kilos"
$origen    = "Provincia de Limarí"
$kgUnid    = 25
$clasif    = "Hortaliza"

# New rows to insert (fecha serial, variedad, calidad, volumen, min, max, promedio, precioKg)
$newRows = @(
    @{ Row = 274; Fecha = 44769; Variedad = "Americana (o)"; Calidad = "Primera"; Volumen = 160; Min = 46000; Max = 48000; Prom = 47000; PrecioKg = 1880 },
    @{ Row = 275; Fecha = 44769; Variedad = "Americana (o)"; Calidad = "Segunda"; Volumen = 120; Min = 36000; Max = 38000; Prom = 37000; PrecioKg = 1480 },
    @{ Row = 276; Fecha = 44769; Variedad = "Inferno";       Calidad = "Primera"; Volumen = 160; Min = 28000; Max = 30000; Prom = 29000; PrecioKg = 1160 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PrecioKg
    $ws.Cells.Item($row, 17).Value = $kgUnid
    $ws.Cells.Item($row, 18).Value = $clasif
}
